$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap rows 136 and 137 (columns B:AC), keep column A (index) fixed
$ws.Range("B136").Value = 6989630
$ws.Range("B137").Value = 6989629
$ws.Range("C136").Value = "Serbia Prva Liga"
$ws.Range("C137").Value = "Serbia Prva Liga"
$ws.Range("D136").Value = "Serbia Prva Liga"
$ws.Range("D137").Value = "Serbia Prva Liga"
$ws.Range("E136").Value = 45256.375
$ws.Range("E137").Value = 45256.375
$ws.Range("F136").Value = "OFK Belgrade"
$ws.Range("F137").Value = "Smederevo"
$ws.Range("G136").Value = "Jedinstvo UB"
$ws.Range("G137").Value = "OFK Vrsac"
$ws.Range("H136").Value = 3
$ws.Range("H137").Value = 1
$ws.Range("I136").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J136").Value = "H"
$ws.Range("J137").Value = "H"
$ws.Range("K136").Value = 1.444
$ws.Range("K137").Value = 2.1
$ws.Range("L136").Value = 3.75
$ws.Range("L137").Value = 2.875
$ws.Range("M136").Value = 6.5
$ws.Range("M137").Value = 3.5
$ws.Range("N136").Value = 1.533
$ws.Range("N137").Value = 1.75
$ws.Range("O136").Value = 3.8
$ws.Range("O137").Value = 3.2
$ws.Range("P136").Value = 6
$ws.Range("P137").Value = 4.75
$ws.Range("Q136").Value = -1.25
$ws.Range("Q137").Value = -0.75
$ws.Range("R136").Value = 2.075
$ws.Range("R137").Value = 1.95
$ws.Range("S136").Value = 1.725
$ws.Range("S137").Value = 1.75
$ws.Range("T136").Value = 2.5
$ws.Range("T137").Value = 2
$ws.Range("U136").Value = 1.9
$ws.Range("U137").Value = 2.025
$ws.Range("V136").Value = 1.9
$ws.Range("V137").Value = 1.775
$ws.Range("W136").Value = 0.5329999999999999
$ws.Range("W137").Value = 0.75
$ws.Range("X136").Value = -1
$ws.Range("X137").Value = -1
$ws.Range("Y136").Value = -1
$ws.Range("Y137").Value = -1
$ws.Range("Z136").Value = 1.075
$ws.Range("Z137").Value = 0.475
$ws.Range("AA136").Value = -1
$ws.Range("AA137").Value = -0.5
$ws.Range("AB136").Value = 0.8999999999999999
$ws.Range("AB137").Value = -1
$ws.Range("AC136").Value = -1
$ws.Range("AC137").Value = 0.7749999999999999

# Swap rows 144 and 147 (columns B:AC), keep column A (index) fixed
$ws.Range("B144").Value = 6989529
$ws.Range("B147").Value = 7019002
$ws.Range("C144").Value = "Serbia Prva Liga"
$ws.Range("C147").Value = "Serbia Prva Liga"
$ws.Range("D144").Value = "Serbia Prva Liga"
$ws.Range("D147").Value = "Serbia Prva Liga"
$ws.Range("E144").Value = 45262.375
$ws.Range("E147").Value = 45262.375
$ws.Range("F144").Value = "FK Radnicki Beograd"
$ws.Range("F147").Value = "Metalac Gornji"
$ws.Range("G144").Value = "OFK Belgrade"
$ws.Range("G147").Value = "FK Kolubara"
$ws.Range("H144").Value = 1
$ws.Range("H147").Value = 2
$ws.Range("I144").Value = 2
$ws.Range("I147").Value = 2
$ws.Range("J144").Value = "A"
$ws.Range("J147").Value = "D"
$ws.Range("K144").Value = 4.75
$ws.Range("K147").Value = 2.5
$ws.Range("L144").Value = 3.75
$ws.Range("L147").Value = 3.2
$ws.Range("M144").Value = 1.571
$ws.Range("M147").Value = 2.5
$ws.Range("N144").Value = 4.75
$ws.Range("N147").Value = 1.75
$ws.Range("O144").Value = 3.75
$ws.Range("O147").Value = 3.1
$ws.Range("P144").Value = 1.571
$ws.Range("P147").Value = 4.5
$ws.Range("Q144").Value = 1
$ws.Range("Q147").Value = -0.5
$ws.Range("R144").Value = 1.8
$ws.Range("R147").Value = 1.8
$ws.Range("S144").Value = 2
$ws.Range("S147").Value = 2
$ws.Range("T144").Value = 2.5
$ws.Range("T147").Value = 2
$ws.Range("U144").Value = 1.85
$ws.Range("U147").Value = 1.9
$ws.Range("V144").Value = 1.95
$ws.Range("V147").Value = 1.9
$ws.Range("W144").Value = -1
$ws.Range("W147").Value = -1
$ws.Range("X144").Value = -1
$ws.Range("X147").Value = 2.1
$ws.Range("Y144").Value = 0.571
$ws.Range("Y147").Value = -1
$ws.Range("Z144").Value = 0
$ws.Range("Z147").Value = -1
$ws.Range("AA144").Value = -0
$ws.Range("AA147").Value = 1
$ws.Range("AB144").Value = 0.8500000000000001
$ws.Range("AB147").Value = 0.8999999999999999
$ws.Range("AC144").Value = -1
$ws.Range("AC147").Value = -1

# Swap rows 145 and 146 (columns B:AC), keep column A (index) fixed
$ws.Range("B145").Value = 6989631
$ws.Range("B146").Value = 6989700
$ws.Range("C145").Value = "Serbia Prva Liga"
$ws.Range("C146").Value = "Serbia Prva Liga"
$ws.Range("D145").Value = "Serbia Prva Liga"
$ws.Range("D146").Value = "Serbia Prva Liga"
$ws.Range("E145").Value = 45262.375
$ws.Range("E146").Value = 45262.375
$ws.Range("F145").Value = "Jedinstvo UB"
$ws.Range("F146").Value = "Sloboda Uzice"
$ws.Range("G145").Value = "FK Tekstilac Odzaci"
$ws.Range("G146").Value = "Radnicki Sremska Mitrovica"
$ws.Range("H145").Value = 1
$ws.Range("H146").Value = 0
$ws.Range("I145").Value = 1
$ws.Range("I146").Value = 1
$ws.Range("J145").Value = "D"
$ws.Range("J146").Value = "A"
$ws.Range("K145").Value = 2.25
$ws.Range("K146").Value = 2.625
$ws.Range("L145").Value = 3
$ws.Range("L146").Value = 3
$ws.Range("M145").Value = 3
$ws.Range("M146").Value = 2.5
$ws.Range("N145").Value = 2.25
$ws.Range("N146").Value = 2.375
$ws.Range("O145").Value = 3
$ws.Range("O146").Value = 3.1
$ws.Range("P145").Value = 3
$ws.Range("P146").Value = 2.75
$ws.Range("Q145").Value = -0.25
$ws.Range("Q146").Value = 0
$ws.Range("R145").Value = 2
$ws.Range("R146").Value = 1.75
$ws.Range("S145").Value = 1.8
$ws.Range("S146").Value = 2.05
$ws.Range("T145").Value = 2.25
$ws.Range("T146").Value = 2
$ws.Range("U145").Value = 2
$ws.Range("U146").Value = 2.025
$ws.Range("V145").Value = 1.8
$ws.Range("V146").Value = 1.775
$ws.Range("W145").Value = -1
$ws.Range("W146").Value = -1
$ws.Range("X145").Value = 2
$ws.Range("X146").Value = -1
$ws.Range("Y145").Value = -1
$ws.Range("Y146").Value = 1.75
$ws.Range("Z145").Value = -0.5
$ws.Range("Z146").Value = -1
$ws.Range("AA145").Value = 0.4
$ws.Range("AA146").Value = 1.05
$ws.Range("AB145").Value = -0.5
$ws.Range("AB146").Value = -1
$ws.Range("AC145").Value = 0.4
$ws.Range("AC146").Value = 0.7749999999999999

# Swap rows 171 and 172 (columns B:AC), keep column A (index) fixed
$ws.Range("B171").Value = 7901796
$ws.Range("B172").Value = 7901795
$ws.Range("C171").Value = "Serbia Prva Liga"
$ws.Range("C172").Value = "Serbia Prva Liga"
$ws.Range("D171").Value = "Serbia Prva Liga"
$ws.Range("D172").Value = "Serbia Prva Liga"
$ws.Range("E171").Value = 45353.375
$ws.Range("E172").Value = 45353.375
$ws.Range("F171").Value = "OFK Belgrade"
$ws.Range("F172").Value = "Radnicki Sremska Mitrovica"
$ws.Range("G171").Value = "FK Graficar Beograd"
$ws.Range("G172").Value = "FK Radnicki Beograd"
$ws.Range("H171").Value = 0
$ws.Range("H172").Value = 1
$ws.Range("I171").Value = 2
$ws.Range("I172").Value = 0
$ws.Range("J171").Value = "A"
$ws.Range("J172").Value = "H"
$ws.Range("K171").Value = 1.4
$ws.Range("K172").Value = 1.666
$ws.Range("L171").Value = 3.6
$ws.Range("L172").Value = 3.4
$ws.Range("M171").Value = 8
$ws.Range("M172").Value = 4.5
$ws.Range("N171").Value = 1.65
$ws.Range("N172").Value = 1.444
$ws.Range("O171").Value = 3.5
$ws.Range("O172").Value = 3.75
$ws.Range("P171").Value = 4.75
$ws.Range("P172").Value = 7
$ws.Range("Q171").Value = -0.75
$ws.Range("Q172").Value = -1
$ws.Range("R171").Value = 1.825
$ws.Range("R172").Value = 1.725
$ws.Range("S171").Value = 1.975
$ws.Range("S172").Value = 1.975
$ws.Range("T171").Value = 2.5
$ws.Range("T172").Value = 2.25
$ws.Range("U171").Value = 1.975
$ws.Range("U172").Value = 1.9
$ws.Range("V171").Value = 1.825
$ws.Range("V172").Value = 1.9
$ws.Range("W171").Value = -1
$ws.Range("W172").Value = 0.444
$ws.Range("X171").Value = -1
$ws.Range("X172").Value = -1
$ws.Range("Y171").Value = 3.75
$ws.Range("Y172").Value = -1
$ws.Range("Z171").Value = -1
$ws.Range("Z172").Value = 0
$ws.Range("AA171").Value = 0.9750000000000001
$ws.Range("AA172").Value = -0
$ws.Range("AB171").Value = -1
$ws.Range("AB172").Value = -1
$ws.Range("AC171").Value = 0.825
$ws.Range("AC172").Value = 0.8999999999999999

# Swap rows 178 and 179 (columns B:AC), keep column A (index) fixed
$ws.Range("B178").Value = 6989710
$ws.Range("B179").Value = 6989640
$ws.Range("C178").Value = "Serbia Prva Liga"
$ws.Range("C179").Value = "Serbia Prva Liga"
$ws.Range("D178").Value = "Serbia Prva Liga"
$ws.Range("D179").Value = "Serbia Prva Liga"
$ws.Range("E178").Value = 45357.41666666666
$ws.Range("E179").Value = 45357.41666666666
$ws.Range("F178").Value = "OFK Belgrade"
$ws.Range("F179").Value = "FK Graficar Beograd"
$ws.Range("G178").Value = "FK Indija"
$ws.Range("G179").Value = "FK Tekstilac Odzaci"
$ws.Range("H178").Value = 0
$ws.Range("H179").Value = 2
$ws.Range("I178").Value = 0
$ws.Range("I179").Value = 2
$ws.Range("J178").Value = "D"
$ws.Range("J179").Value = "D"
$ws.Range("K178").Value = 1.571
$ws.Range("K179").Value = 1.5
$ws.Range("L178").Value = 3.25
$ws.Range("L179").Value = 3.4
$ws.Range("M178").Value = 6
$ws.Range("M179").Value = 6.5
$ws.Range("N178").Value = 1.75
$ws.Range("N179").Value = 1.85
$ws.Range("O178").Value = 3
$ws.Range("O179").Value = 3.2
$ws.Range("P178").Value = 4.75
$ws.Range("P179").Value = 3.75
$ws.Range("Q178").Value = -0.5
$ws.Range("Q179").Value = -0.5
$ws.Range("R178").Value = 1.8
$ws.Range("R179").Value = 1.95
$ws.Range("S178").Value = 2
$ws.Range("S179").Value = 1.85
$ws.Range("T178").Value = 2.25
$ws.Range("T179").Value = 2
$ws.Range("U178").Value = 1.85
$ws.Range("U179").Value = 1.775
$ws.Range("V178").Value = 1.95
$ws.Range("V179").Value = 2.025
$ws.Range("W178").Value = -1
$ws.Range("W179").Value = -1
$ws.Range("X178").Value = 2
$ws.Range("X179").Value = 2.2
$ws.Range("Y178").Value = -1
$ws.Range("Y179").Value = -1
$ws.Range("Z178").Value = -1
$ws.Range("Z179").Value = -1
$ws.Range("AA178").Value = 1
$ws.Range("AA179").Value = 0.8500000000000001
$ws.Range("AB178").Value = -1
$ws.Range("AB179").Value = 0.7749999999999999
$ws.Range("AC178").Value = 0.95
$ws.Range("AC179").Value = -1

# Swap rows 213 and 214 (columns B:AC), keep column A (index) fixed
$ws.Range("B213").Value = 6989720
$ws.Range("B214").Value = 6989719
$ws.Range("C213").Value = "Serbia Prva Liga"
$ws.Range("C214").Value = "Serbia Prva Liga"
$ws.Range("D213").Value = "Serbia Prva Liga"
$ws.Range("D214").Value = "Serbia Prva Liga"
$ws.Range("E213").Value = 45381.45833333334
$ws.Range("E214").Value = 45381.45833333334
$ws.Range("F213").Value = "OFK Belgrade"
$ws.Range("F214").Value = "FK Tekstilac Odzaci"
$ws.Range("G213").Value = "Smederevo"
$ws.Range("G214").Value = "Radnicki Sremska Mitrovica"
$ws.Range("H213").Value = 1
$ws.Range("H214").Value = 5
$ws.Range("I213").Value = 2
$ws.Range("I214").Value = 2
$ws.Range("J213").Value = "A"
$ws.Range("J214").Value = "H"
$ws.Range("K213").Value = 1.444
$ws.Range("K214").Value = 1.833
$ws.Range("L213").Value = 4.2
$ws.Range("L214").Value = 3
$ws.Range("M213").Value = 5.5
$ws.Range("M214").Value = 4.2
$ws.Range("N213").Value = 1.444
$ws.Range("N214").Value = 1.7
$ws.Range("O213").Value = 4
$ws.Range("O214").Value = 3
$ws.Range("P213").Value = 6
$ws.Range("P214").Value = 5.25
$ws.Range("Q213").Value = -1
$ws.Range("Q214").Value = -0.75
$ws.Range("R213").Value = 1.725
$ws.Range("R214").Value = 2
$ws.Range("S213").Value = 1.975
$ws.Range("S214").Value = 1.8
$ws.Range("T213").Value = 2.25
$ws.Range("T214").Value = 1.75
$ws.Range("U213").Value = 1.8
$ws.Range("U214").Value = 1.9
$ws.Range("V213").Value = 2
$ws.Range("V214").Value = 1.9
$ws.Range("W213").Value = -1
$ws.Range("W214").Value = 0.7
$ws.Range("X213").Value = -1
$ws.Range("X214").Value = -1
$ws.Range("Y213").Value = 5
$ws.Range("Y214").Value = -1
$ws.Range("Z213").Value = -1
$ws.Range("Z214").Value = 1
$ws.Range("AA213").Value = 0.9750000000000001
$ws.Range("AA214").Value = -1
$ws.Range("AB213").Value = 0.8
$ws.Range("AB214").Value = 0.8999999999999999
$ws.Range("AC213").Value = -1
$ws.Range("AC214").Value = -1
